$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1073.3334
$ws.Range("I98").Value = 888.5143
$ws.Range("J98").Value = 1997.4286
$ws.Range("K98").Value = 888.5143
$ws.Range("L98").Value = 1997.4286
$ws.Range("M98").Value = 609.4857
$ws.Range("N98").Value = -4993.4286

$ws.Range("H122").Value = 1073.3334
$ws.Range("I122").Value = 888.5143
$ws.Range("J122").Value = 1997.4286
$ws.Range("K122").Value = 2665.5429
$ws.Range("L122").Value = 5992.2858
$ws.Range("M122").Value = -215.5429000000004
$ws.Range("N122").Value = -10892.2858

$ws.Range("H132").Value = 5381.074
$ws.Range("I132").Value = 5268.8945
$ws.Range("J132").Value = 5647.5
$ws.Range("K132").Value = 15806.6835
$ws.Range("L132").Value = 16942.5
$ws.Range("M132").Value = -13276.6835
$ws.Range("N132").Value = -22002.5

$ws.Range("H137").Value = 15983.888
$ws.Range("I137").Value = 956.86957
$ws.Range("J137").Value = 43633.6
$ws.Range("K137").Value = 2870.60871
$ws.Range("L137").Value = 130900.8
$ws.Range("M137").Value = -320.60871
$ws.Range("N137").Value = -136000.8

$ws.Range("H138").Value = 1505.0613
$ws.Range("I138").Value = 890.1212
$ws.Range("J138").Value = 2773.375
$ws.Range("K138").Value = 2670.3636
$ws.Range("L138").Value = 8320.125
$ws.Range("M138").Value = 2469.6364
$ws.Range("N138").Value = -18600.125

$ws.Range("H141").Value = 744.5577
$ws.Range("I141").Value = 511.48935
$ws.Range("J141").Value = 2935.4
$ws.Range("K141").Value = 1534.46805
$ws.Range("L141").Value = 8806.200000000001
$ws.Range("M141").Value = 3645.53195
$ws.Range("N141").Value = -19166.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1379.28
$ws.Range("I32").Value = 1229.4584
$ws.Range("J32").Value = 4975
$ws.Range("K32").Value = 1229.4584
$ws.Range("L32").Value = 4975
$ws.Range("M32").Value = -942.4584
$ws.Range("N32").Value = -5549

$ws.Range("H61").Value = 1429.16
$ws.Range("I61").Value = 1152.659
$ws.Range("J61").Value = 3456.8333
$ws.Range("K61").Value = 1152.659
$ws.Range("L61").Value = 3456.8333
$ws.Range("M61").Value = -940.6590000000001
$ws.Range("N61").Value = -3880.8333

$ws.Range("H136").Value = 1429.16
$ws.Range("I136").Value = 1152.659
$ws.Range("J136").Value = 3456.8333
$ws.Range("K136").Value = 3457.977
$ws.Range("L136").Value = 10370.4999
$ws.Range("M136").Value = -907.9770000000003
$ws.Range("N136").Value = -15470.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2006.6111
$ws.Range("I134").Value = 1739.5483
$ws.Range("J134").Value = 3662.4
$ws.Range("K134").Value = 5218.644899999999
$ws.Range("L134").Value = 10987.2
$ws.Range("M134").Value = -2683.644899999999
$ws.Range("N134").Value = -16057.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2321.8406
$ws.Range("I31").Value = 1194.5834
$ws.Range("J31").Value = 4898.4287
$ws.Range("K31").Value = 1194.5834
$ws.Range("L31").Value = 4898.4287
$ws.Range("M31").Value = -899.5834
$ws.Range("N31").Value = -5488.4287

$ws.Range("H34").Value = 2321.8406
$ws.Range("I34").Value = 1194.5834
$ws.Range("J34").Value = 4898.4287
$ws.Range("K34").Value = 1194.5834
$ws.Range("L34").Value = 4898.4287
$ws.Range("M34").Value = -992.5834
$ws.Range("N34").Value = -5302.4287

$ws.Range("H132").Value = 2318.7
$ws.Range("I132").Value = 2263.375
$ws.Range("J132").Value = 2540
$ws.Range("K132").Value = 6790.125
$ws.Range("L132").Value = 7620
$ws.Range("M132").Value = -4260.125
$ws.Range("N132").Value = -12680

$ws.Range("H134").Value = 1860.7903
$ws.Range("I134").Value = 1239.3864
$ws.Range("J134").Value = 3379.7778
$ws.Range("K134").Value = 3718.1592
$ws.Range("L134").Value = 10139.3334
$ws.Range("M134").Value = -1183.1592
$ws.Range("N134").Value = -15209.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 216.72223
$ws.Range("J107").Value = 175.61539
$ws.Range("L107").Value = 526.84617
$ws.Range("N107").Value = -4366.84617

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 854.5
$ws.Range("I97").Value = 502.14285
$ws.Range("J97").Value = 1676.6666
$ws.Range("K97").Value = 502.14285
$ws.Range("L97").Value = 1676.6666
$ws.Range("M97").Value = -6.14285000000001
$ws.Range("N97").Value = -2668.6666

$ws.Range("H122").Value = 1473.3462
$ws.Range("I122").Value = 1350.0555
$ws.Range("J122").Value = 1750.75
$ws.Range("K122").Value = 4050.1665
$ws.Range("L122").Value = 5252.25
$ws.Range("M122").Value = -1600.1665
$ws.Range("N122").Value = -10152.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 10000
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("N96").Value = -15492

$ws.Range("H100").Value = 500000000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 500000000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 500000000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -500001082

$ws.Range("H132").Value = 1701.42
$ws.Range("I132").Value = 1718.6882
$ws.Range("J132").Value = 1472
$ws.Range("K132").Value = 5156.0646
$ws.Range("L132").Value = 4416
$ws.Range("M132").Value = -2626.0646
$ws.Range("N132").Value = -9476

$ws.Range("H136").Value = 2034.4492
$ws.Range("I136").Value = 1566.1273
$ws.Range("J136").Value = 3874.2856
$ws.Range("K136").Value = 4698.3819
$ws.Range("L136").Value = 11622.8568
$ws.Range("M136").Value = -2148.3819
$ws.Range("N136").Value = -16722.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1820.5897
$ws.Range("I122").Value = 1589.8064
$ws.Range("J122").Value = 2714.875
$ws.Range("K122").Value = 4769.4192
$ws.Range("L122").Value = 8144.625
$ws.Range("M122").Value = -2319.4192
$ws.Range("N122").Value = -13044.625

$ws.Range("H126").Value = 2689.611
$ws.Range("I126").Value = 2006.5333
$ws.Range("J126").Value = 6105
$ws.Range("K126").Value = 6019.5999
$ws.Range("L126").Value = 18315
$ws.Range("M126").Value = -3549.5999
$ws.Range("N126").Value = -23255

$ws.Range("H132").Value = 1113.0834
$ws.Range("I132").Value = 654.26984
$ws.Range("J132").Value = 2489.524
$ws.Range("K132").Value = 1962.80952
$ws.Range("L132").Value = 7468.572
$ws.Range("M132").Value = 567.1904799999998
$ws.Range("N132").Value = -12528.572
